$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2553.2969
$ws.Range("I15").Value = 2553.2969
$ws.Range("K15").Value = 7659.8907
$ws.Range("M15").Value = -7490.8907
$ws.Range("H112").Value = 1845.4166
$ws.Range("J112").Value = 1903.9131
$ws.Range("L112").Value = 5711.7393
$ws.Range("N112").Value = -7927.7393
$ws.Range("H135").Value = 1028.5333
$ws.Range("I135").Value = 472.45
$ws.Range("K135").Value = 4252.05
$ws.Range("M135").Value = -1717.05
$ws.Range("H137").Value = 3532
$ws.Range("I137").Value = 1215.4117
$ws.Range("J137").Value = 4469.6665
$ws.Range("K137").Value = 3646.2351
$ws.Range("L137").Value = 13408.9995
$ws.Range("M137").Value = -1096.2351
$ws.Range("N137").Value = -18508.9995

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1365.7812
$ws.Range("I2").Value = 1233.8077
$ws.Range("J2").Value = 1937.6666
$ws.Range("K2").Value = 1233.8077
$ws.Range("L2").Value = 1937.6666
$ws.Range("M2").Value = -1120.8077
$ws.Range("N2").Value = -2163.6666
$ws.Range("H32").Value = 5524.357
$ws.Range("I32").Value = 6181.794
$ws.Range("J32").Value = 2730.25
$ws.Range("K32").Value = 6181.794
$ws.Range("L32").Value = 2730.25
$ws.Range("M32").Value = -5894.794
$ws.Range("N32").Value = -3304.25
$ws.Range("H45").Value = 1648.5
$ws.Range("I45").Value = 1116.4445
$ws.Range("J45").Value = 2332.5715
$ws.Range("K45").Value = 1116.4445
$ws.Range("L45").Value = 2332.5715
$ws.Range("M45").Value = -739.4445000000001
$ws.Range("N45").Value = -3086.5715
$ws.Range("H61").Value = 25001706
$ws.Range("I61").Value = 31251474
$ws.Range("J61").Value = 2627.75
$ws.Range("K61").Value = 31251474
$ws.Range("L61").Value = 2627.75
$ws.Range("M61").Value = -31251262
$ws.Range("N61").Value = -3051.75
$ws.Range("H74").Value = 4008.6365
$ws.Range("I74").Value = 764.2917
$ws.Range("J74").Value = 12660.223
$ws.Range("K74").Value = 764.2917
$ws.Range("L74").Value = 12660.223
$ws.Range("M74").Value = 109.7083
$ws.Range("N74").Value = -14408.223
$ws.Range("H77").Value = 4008.6365
$ws.Range("I77").Value = 764.2917
$ws.Range("J77").Value = 12660.223
$ws.Range("K77").Value = 3821.4585
$ws.Range("L77").Value = 63301.115
$ws.Range("M77").Value = 546.5415000000003
$ws.Range("N77").Value = -72037.11499999999
$ws.Range("H110").Value = 1128.2727
$ws.Range("I110").Value = 1114.0358
$ws.Range("J110").Value = 1208
$ws.Range("K110").Value = 1114.0358
$ws.Range("L110").Value = 1208
$ws.Range("M110").Value = 930.9641999999999
$ws.Range("N110").Value = -5298
$ws.Range("H116").Value = 1365.7812
$ws.Range("I116").Value = 1233.8077
$ws.Range("J116").Value = 1937.6666
$ws.Range("K116").Value = 1233.8077
$ws.Range("L116").Value = 1937.6666
$ws.Range("M116").Value = 1060.1923
$ws.Range("N116").Value = -6525.6666
$ws.Range("H122").Value = 1424.8889
$ws.Range("I122").Value = 1424.8889
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4274.6667
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1824.6667
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 25001706
$ws.Range("I136").Value = 31251474
$ws.Range("J136").Value = 2627.75
$ws.Range("K136").Value = 93754422
$ws.Range("L136").Value = 7883.25
$ws.Range("M136").Value = -93751872
$ws.Range("N136").Value = -12983.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1365.7812
$ws.Range("I3").Value = 1233.8077
$ws.Range("J3").Value = 1937.6666
$ws.Range("K3").Value = 1233.8077
$ws.Range("L3").Value = 1937.6666
$ws.Range("M3").Value = -1119.8077
$ws.Range("N3").Value = -2165.6666
$ws.Range("H35").Value = 30000
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H99").Value = 1017.6875
$ws.Range("I99").Value = 872.2222
$ws.Range("J99").Value = 1204.7142
$ws.Range("K99").Value = 872.2222
$ws.Range("L99").Value = 1204.7142
$ws.Range("M99").Value = 625.7778
$ws.Range("N99").Value = -4200.7142
$ws.Range("H107").Value = 1405.7931
$ws.Range("I107").Value = 1040.6316
$ws.Range("J107").Value = 2099.6
$ws.Range("K107").Value = 1040.6316
$ws.Range("L107").Value = 2099.6
$ws.Range("M107").Value = 879.3684000000001
$ws.Range("N107").Value = -5939.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2834.0625
$ws.Range("I16").Value = 1919.3334
$ws.Range("J16").Value = 5578.25
$ws.Range("K16").Value = 1919.3334
$ws.Range("L16").Value = 5578.25
$ws.Range("M16").Value = -1632.3334
$ws.Range("N16").Value = -6152.25
$ws.Range("H31").Value = 7899.6562
$ws.Range("I31").Value = 2249.1
$ws.Range("J31").Value = 10468.091
$ws.Range("K31").Value = 2249.1
$ws.Range("L31").Value = 10468.091
$ws.Range("M31").Value = -1954.1
$ws.Range("N31").Value = -11058.091
$ws.Range("H34").Value = 7899.6562
$ws.Range("I34").Value = 2249.1
$ws.Range("J34").Value = 10468.091
$ws.Range("K34").Value = 2249.1
$ws.Range("L34").Value = 10468.091
$ws.Range("M34").Value = -2047.1
$ws.Range("N34").Value = -10872.091
$ws.Range("H113").Value = 2834.0625
$ws.Range("I113").Value = 1919.3334
$ws.Range("J113").Value = 5578.25
$ws.Range("K113").Value = 1919.3334
$ws.Range("L113").Value = 5578.25
$ws.Range("M113").Value = 250.6666
$ws.Range("N113").Value = -9918.25
$ws.Range("H122").Value = 863.875
$ws.Range("I122").Value = 690.2778
$ws.Range("J122").Value = 1384.6666
$ws.Range("K122").Value = 2070.8334
$ws.Range("L122").Value = 4153.9998
$ws.Range("M122").Value = 379.1666
$ws.Range("N122").Value = -9053.9998
$ws.Range("H134").Value = 2832
$ws.Range("I134").Value = 2994.6667
$ws.Range("J134").Value = 2588
$ws.Range("K134").Value = 8984.000100000001
$ws.Range("L134").Value = 7764
$ws.Range("M134").Value = -6449.000100000001
$ws.Range("N134").Value = -12834

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 219
$ws.Range("I107").Value = 207.42857
$ws.Range("K107").Value = 207.42857
$ws.Range("M107").Value = 1712.57143
$ws.Range("H113").Value = 3492.2144
$ws.Range("I113").Value = 3346.375
$ws.Range("J113").Value = 3686.6667
$ws.Range("K113").Value = 3346.375
$ws.Range("L113").Value = 3686.6667
$ws.Range("M113").Value = -1176.375
$ws.Range("N113").Value = -8026.6667
$ws.Range("H122").Value = 1261.12
$ws.Range("I122").Value = 1271.4
$ws.Range("J122").Value = 1220
$ws.Range("K122").Value = 3814.2
$ws.Range("L122").Value = 3660
$ws.Range("M122").Value = -1364.2
$ws.Range("N122").Value = -8560

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 48000
$ws.Range("J63").Value = 48000
$ws.Range("L63").Value = 48000
$ws.Range("N63").Value = -49498
$ws.Range("H66").Value = 48000
$ws.Range("J66").Value = 48000
$ws.Range("L66").Value = 144000
$ws.Range("N66").Value = -151488
$ws.Range("H98").Value = 48000
$ws.Range("J98").Value = 48000
$ws.Range("L98").Value = 48000
$ws.Range("N98").Value = -53990

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 820.5714
$ws.Range("I96").Value = 595
$ws.Range("J96").Value = 910.8
$ws.Range("K96").Value = 595
$ws.Range("L96").Value = 910.8
$ws.Range("M96").Value = 778
$ws.Range("N96").Value = -3656.8
$ws.Range("H132").Value = 3873.7192
$ws.Range("I132").Value = 4291.4
$ws.Range("J132").Value = 2890.9412
$ws.Range("K132").Value = 12874.2
$ws.Range("L132").Value = 8672.8236
$ws.Range("N132").Value = -13732.8236
$ws.Range("M132").Value = -10344.2
